$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet
$ws.Name = "MV_Testdaten_Iterations"

# Expand the table to a 3rd column (new column becomes "Zulassungsland",
# a copy of what used to be column A) before we repurpose column A.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:C4"))

# Copy the old "Zulassungsland" column (A) into the new column (C).
$ws.Range("C1").Value = "Zulassungsland"
$ws.Range("C2").Value = "DE"
$ws.Range("C3").Value = "PL"
$ws.Range("C4").Value = "DE"

# Repurpose column A into the new "Buchungssprache" (booking language) column.
$ws.Range("A1").Value = "Buchungssprache"
$ws.Range("A2").Value = "Deutsch"
$ws.Range("A3").Value = "Polski"
$ws.Range("A4").Value = "Deutsch"

# Column B ("Kennzeichen") stays as-is.

# Match column widths as closely as this engine's width model allows:
# column A got wider (new "Buchungssprache" header), column C takes on the
# width the old "Zulassungsland" column (A) used to have.
$ws.Columns("A").ColumnWidth = 18.5
$ws.Columns("C").ColumnWidth = 16

# Restore the selected cell as recorded in the saved view.
[void]$ws.Range("C31").Select()
